# Add a "Fade" slide transition (medium speed, 0.7s / 700ms duration) to
# every slide in the deck. This mirrors PowerPoint authoring the Fade
# transition (Transitions ribbon -> Fade, Duration 00.70) on all slides.
#
# PpEntryEffect.ppEffectFadeSmoothly = 1793 (Fade transition "Smoothly").
# PpTransitionSpeed.ppTransitionSpeedMedium = 2.
# Duration is expressed in seconds on the object model (0.7s -> 700ms
# in the underlying OOXML p14:dur attribute).

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $t = $slide.SlideShowTransition

    # Order matters: set the entry effect (and implied duration) first,
    # then the explicit duration, and set Speed last so the "medium"
    # speed sticks instead of being overwritten by the effect/duration
    # assignment.
    $t.EntryEffect = 1793
    $t.Duration = 0.7
    $t.Speed = 2
}
